$wb = $excel.ActiveWorkbook

# --- Summary ---
$ws = $wb.Worksheets.Item("Summary")
$ws.Range("B4").Value = "inf"
$ws.Range("B7").Value = 10414278.58937073
$ws.Range("B8").Value = 26920418.46524543
$ws.Range("B10").Value = 1992691.030333921

# --- Fed-in Capacity ---
$ws = $wb.Worksheets.Item("Fed-in Capacity")
$ws.Range("K2").Value = 214.587604768856
$ws.Range("L2").Value = 229.8722545957376
$ws.Range("M2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 225.4521708613878
$ws.Range("Q2").Value = 216.7480476275882
$ws.Range("K3").Value = 134.395403
$ws.Range("N3").Value = 128.05816928125
$ws.Range("O3").Value = 139.0313383333333
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("M4").Value = 135.4526393489149
$ws.Range("O4").Value = 0
$ws.Range("K5").Value = 214.587604768856
$ws.Range("N6").Value = 0
$ws.Range("P6").Value = 130.625047228972
$ws.Range("Q6").Value = 0
$ws.Range("O7").Value = 134.9951249905467
$ws.Range("K8").Value = 214.587604768856
$ws.Range("M8").Value = 224.5875773965909
$ws.Range("N8").Value = 0
$ws.Range("Q8").Value = 216.7480476275882
$ws.Range("L9").Value = 0
$ws.Range("P9").Value = 0
$ws.Range("Q9").Value = 136.482229733871
$ws.Range("O10").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").Value = 224.5875773965909
$ws.Range("Q11").Value = 216.7480476275882
$ws.Range("M12").Value = 138.5806830739679
$ws.Range("O12").Value = 139.0313383333333
$ws.Range("Q12").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = 128.05816928125
$ws.Range("O15").Value = 0
$ws.Range("Q15").Value = 136.482229733871
$ws.Range("N16").Value = 0
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = 134.2848039479189
$ws.Range("K17").Value = 214.587604768856
$ws.Range("L17").Value = 229.8722545957376
$ws.Range("Q17").Value = 216.7480476275882
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = 138.5806830739679
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = 130.625047228972
$ws.Range("Q18").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 0
$ws.Range("O19").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("N20").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("O21").Value = 139.0313383333333
$ws.Range("P21").Value = 130.625047228972
$ws.Range("Q21").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 0
$ws.Range("P22").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").Value = 224.5875773965909
$ws.Range("O23").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("P24").Value = 130.625047228972
$ws.Range("Q24").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 0
$ws.Range("P25").Value = 0
$ws.Range("L26").Value = 229.8722545957376
$ws.Range("O26").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("P27").Value = 0
$ws.Range("O28").Value = 0
$ws.Range("L29").Value = 229.8722545957376
$ws.Range("M29").Value = 224.5875773965909
$ws.Range("N29").Value = 223.6777370066762
$ws.Range("Q29").Value = 216.7480476275882
$ws.Range("M30").Value = 138.5806830739679
$ws.Range("N30").Value = 0
$ws.Range("P30").Value = 130.625047228972
$ws.Range("L31").Value = 131.5125593742073
$ws.Range("N31").Value = 124.4934058536024
$ws.Range("O31").Value = 0
$ws.Range("P31").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = 0
$ws.Range("O32").Value = 0
$ws.Range("Q32").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 135.0905202853774
$ws.Range("O33").Value = 0
$ws.Range("P33").Value = 130.625047228972
$ws.Range("L34").Value = 131.5125593742073
$ws.Range("M34").Value = 135.4526393489149
$ws.Range("N34").Value = 124.4934058536024
$ws.Range("O34").Value = 134.9951249905467
$ws.Range("L35").Value = 0
$ws.Range("P35").Value = 0
$ws.Range("K36").Value = 134.395403
$ws.Range("M36").Value = 138.5806830739679
$ws.Range("M37").Value = 135.4526393489149
$ws.Range("N37").Value = 0
$ws.Range("O37").Value = 134.9951249905467
$ws.Range("P37").Value = 0
$ws.Range("K38").Value = 214.587604768856
$ws.Range("L38").Value = 229.8722545957376
$ws.Range("M38").Value = 0
$ws.Range("Q38").Value = 216.7480476275882
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 134.395403
$ws.Range("N39").Value = 0
$ws.Range("P39").Value = 0
$ws.Range("Q39").Value = 136.482229733871
$ws.Range("N40").Value = 124.4934058536024
$ws.Range("O40").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").Value = 0
$ws.Range("Q41").Value = 216.7480476275882
$ws.Range("J42").Value = 123.666686
$ws.Range("L42").Value = 135.0905202853774
$ws.Range("M42").Value = 0
$ws.Range("O42").Value = 0
$ws.Range("P42").Value = 0
$ws.Range("N43").Value = 124.4934058536024
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("O44").Value = 0
$ws.Range("P44").Value = 0
$ws.Range("Q44").Value = 216.7480476275882
$ws.Range("J45").Value = 0
$ws.Range("O45").Value = 0
$ws.Range("P45").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = 135.4526393489149
$ws.Range("O46").Value = 0
$ws.Range("P46").Value = 0

# --- Unmet Demand ---
$ws = $wb.Worksheets.Item("Unmet Demand")
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 224.5875773965909
$ws.Range("O2").Value = 224.3457561361446
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 130.625047228972
$ws.Range("Q3").Value = 136.482229733871
$ws.Range("M4").Value = 0
$ws.Range("O4").Value = 134.9951249905467
$ws.Range("K5").Value = 0
$ws.Range("N6").Value = 128.05816928125
$ws.Range("P6").Value = 0
$ws.Range("Q6").Value = 136.482229733871
$ws.Range("O7").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").Value = 0
$ws.Range("N8").Value = 223.6777370066762
$ws.Range("Q8").Value = 0
$ws.Range("L9").Value = 135.0905202853774
$ws.Range("P9").Value = 130.625047228972
$ws.Range("Q9").Value = 0
$ws.Range("O10").Value = 134.9951249905467
$ws.Range("K11").Value = 214.587604768856
$ws.Range("M11").Value = 0
$ws.Range("Q11").Value = 0
$ws.Range("M12").Value = 0
$ws.Range("O12").Value = 0
$ws.Range("Q12").Value = 136.482229733871
$ws.Range("L13").Value = 131.5125593742073
$ws.Range("L14").Value = 229.8722545957376
$ws.Range("J15").Value = 123.666686
$ws.Range("L15").Value = 135.0905202853774
$ws.Range("N15").Value = 0
$ws.Range("O15").Value = 139.0313383333333
$ws.Range("Q15").Value = 0
$ws.Range("N16").Value = 124.4934058536024
$ws.Range("O16").Value = 134.9951249905467
$ws.Range("P16").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("Q17").Value = 0
$ws.Range("K18").Value = 134.395403
$ws.Range("L18").Value = 135.0905202853774
$ws.Range("M18").Value = 0
$ws.Range("O18").Value = 139.0313383333333
$ws.Range("P18").Value = 0
$ws.Range("Q18").Value = 136.482229733871
$ws.Range("L19").Value = 131.5125593742073
$ws.Range("M19").Value = 135.4526393489149
$ws.Range("O19").Value = 134.9951249905467
$ws.Range("K20").Value = 214.587604768856
$ws.Range("N20").Value = 223.6777370066762
$ws.Range("J21").Value = 123.666686
$ws.Range("O21").Value = 0
$ws.Range("P21").Value = 0
$ws.Range("Q21").Value = 136.482229733871
$ws.Range("L22").Value = 131.5125593742073
$ws.Range("M22").Value = 135.4526393489149
$ws.Range("P22").Value = 134.2848039479189
$ws.Range("K23").Value = 214.587604768856
$ws.Range("M23").Value = 0
$ws.Range("O23").Value = 224.3457561361446
$ws.Range("K24").Value = 134.395403
$ws.Range("P24").Value = 0
$ws.Range("Q24").Value = 136.482229733871
$ws.Range("L25").Value = 131.5125593742073
$ws.Range("M25").Value = 135.4526393489149
$ws.Range("P25").Value = 134.2848039479189
$ws.Range("L26").Value = 0
$ws.Range("O26").Value = 224.3457561361446
$ws.Range("L27").Value = 135.0905202853774
$ws.Range("P27").Value = 130.625047228972
$ws.Range("O28").Value = 134.9951249905467
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 0
$ws.Range("N29").Value = 0
$ws.Range("Q29").Value = 0
$ws.Range("M30").Value = 0
$ws.Range("N30").Value = 128.05816928125
$ws.Range("P30").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").Value = 0
$ws.Range("O31").Value = 134.9951249905467
$ws.Range("P31").Value = 134.2848039479189
$ws.Range("L32").Value = 229.8722545957376
$ws.Range("M32").Value = 224.5875773965909
$ws.Range("O32").Value = 224.3457561361446
$ws.Range("Q32").Value = 216.7480476275882
$ws.Range("J33").Value = 123.666686
$ws.Range("L33").Value = 0
$ws.Range("O33").Value = 139.0313383333333
$ws.Range("P33").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = 0
$ws.Range("N34").Value = 0
$ws.Range("O34").Value = 0
$ws.Range("L35").Value = 229.8722545957376
$ws.Range("P35").Value = 225.4521708613878
$ws.Range("K36").Value = 0
$ws.Range("M36").Value = 0
$ws.Range("M37").Value = 0
$ws.Range("N37").Value = 124.4934058536024
$ws.Range("O37").Value = 0
$ws.Range("P37").Value = 134.2848039479189
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 224.5875773965909
$ws.Range("Q38").Value = 0
$ws.Range("J39").Value = 123.666686
$ws.Range("K39").Value = 0
$ws.Range("N39").Value = 128.05816928125
$ws.Range("P39").Value = 130.625047228972
$ws.Range("Q39").Value = 0
$ws.Range("N40").Value = 0
$ws.Range("O40").Value = 134.9951249905467
$ws.Range("L41").Value = 229.8722545957376
$ws.Range("N41").Value = 223.6777370066762
$ws.Range("Q41").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = 138.5806830739679
$ws.Range("O42").Value = 139.0313383333333
$ws.Range("P42").Value = 130.625047228972
$ws.Range("N43").Value = 0
$ws.Range("K44").Value = 214.587604768856
$ws.Range("L44").Value = 229.8722545957376
$ws.Range("O44").Value = 224.3457561361446
$ws.Range("P44").Value = 225.4521708613878
$ws.Range("Q44").Value = 0
$ws.Range("J45").Value = 123.666686
$ws.Range("O45").Value = 139.0313383333333
$ws.Range("P45").Value = 130.625047228972
$ws.Range("L46").Value = 131.5125593742073
$ws.Range("M46").Value = 0
$ws.Range("O46").Value = 134.9951249905467
$ws.Range("P46").Value = 134.2848039479189

# --- Household Surplus ---
$ws = $wb.Worksheets.Item("Household Surplus")
$ws.Range("B2").Value = 210488.8061910179
$ws.Range("B3").Value = 220774.506974527
$ws.Range("B4").Value = 211732.8194477764
$ws.Range("B5").Value = 200369.2909544114
$ws.Range("B6").Value = 216198.2868781213
$ws.Range("B7").Value = 180285.0845643228
$ws.Range("B8").Value = 96738.99602561745
$ws.Range("B9").Value = 83897.15139637531
$ws.Range("B10").Value = 155449.9532625354
$ws.Range("B11").Value = 242367.2892298528
$ws.Range("B12").Value = 93240.28610310853
$ws.Range("B13").Value = 166774.6036979506
$ws.Range("B14").Value = 165444.913715929
$ws.Range("B15").Value = 179855.8546410235
$ws.Range("B16").Value = 97473.70738772406
